$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Quiero una senal de 1 kHz" (TIMER0) block: change Pre-escaler Timer 0 from 8 to 1
$ws.Range("D12").Value = 1

# TIMER2 block: change F(PWM) Hz from 9600 to 2400
$ws.Range("D19").Value = 2400

# Recalculate so the dependent formula cells (C14, B21, C21) refresh with the new inputs
$excel.Calculate()

# Restore the view: scrolled up one row and the active selection moved to C19
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C19").Select()
